$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1325
$ws.Range("I127").Value = 659.8
$ws.Range("J127").Value = 1694.5555
$ws.Range("K127").Value = 1979.4
$ws.Range("L127").Value = 5083.666499999999
$ws.Range("M127").Value = 2980.6
$ws.Range("N127").Value = -15003.6665
$ws.Range("H129").Value = 1571.8572
$ws.Range("I129").Value = 1197.75
$ws.Range("J129").Value = 1634.2084
$ws.Range("K129").Value = 3593.25
$ws.Range("L129").Value = 4902.6252
$ws.Range("M129").Value = 1406.75
$ws.Range("N129").Value = -14902.6252
$ws.Range("H132").Value = 6787.722
$ws.Range("I132").Value = 6288.2905
$ws.Range("J132").Value = 7460.8696
$ws.Range("K132").Value = 18864.8715
$ws.Range("L132").Value = 22382.6088
$ws.Range("M132").Value = -16334.8715
$ws.Range("N132").Value = -27442.6088
$ws.Range("H137").Value = 1435.9259
$ws.Range("I137").Value = 1871.7407
$ws.Range("J137").Value = 1000.1111
$ws.Range("K137").Value = 5615.2221
$ws.Range("L137").Value = 3000.3333
$ws.Range("M137").Value = -3065.2221
$ws.Range("N137").Value = -8100.3333
$ws.Range("H138").Value = 2375.8086
$ws.Range("I138").Value = 1443.8148
$ws.Range("J138").Value = 3634
$ws.Range("K138").Value = 4331.4444
$ws.Range("L138").Value = 10902
$ws.Range("M138").Value = 808.5555999999997
$ws.Range("N138").Value = -21182
$ws.Range("H141").Value = 4861.8
$ws.Range("I141").Value = 1465.1515
$ws.Range("J141").Value = 9956.772000000001
$ws.Range("K141").Value = 4395.4545
$ws.Range("L141").Value = 29870.316
$ws.Range("M141").Value = 784.5455000000002
$ws.Range("N141").Value = -40230.31600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2927.65
$ws.Range("I61").Value = 5142.875
$ws.Range("J61").Value = 1450.8334
$ws.Range("K61").Value = 5142.875
$ws.Range("L61").Value = 1450.8334
$ws.Range("M61").Value = -4930.875
$ws.Range("H69").Value = 40997.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 40997.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 40997.5
$ws.Range("N69").Value = -42495.5
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 40997.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 40997.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 122992.5
$ws.Range("N72").Value = -130480.5
$ws.Range("M72").ClearContents()
$ws.Range("H136").Value = 2927.65
$ws.Range("I136").Value = 5142.875
$ws.Range("J136").Value = 1450.8334
$ws.Range("K136").Value = 15428.625
$ws.Range("L136").Value = 4352.5002
$ws.Range("M136").Value = -12878.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2243.7273
$ws.Range("I20").Value = 1304.5
$ws.Range("J20").Value = 2780.4285
$ws.Range("K20").Value = 1304.5
$ws.Range("L20").Value = 2780.4285
$ws.Range("M20").Value = -1057.5
$ws.Range("N20").Value = -3274.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2060.7407
$ws.Range("I132").Value = 1317.6923
$ws.Range("J132").Value = 2750.7144
$ws.Range("K132").Value = 3953.0769
$ws.Range("L132").Value = 8252.143199999999
$ws.Range("M132").Value = -1423.0769
$ws.Range("N132").Value = -13312.1432
$ws.Range("H134").Value = 1960.6923
$ws.Range("I134").Value = 1131.8125
$ws.Range("J134").Value = 2537.3044
$ws.Range("K134").Value = 3395.4375
$ws.Range("L134").Value = 7611.9132
$ws.Range("M134").Value = -860.4375
$ws.Range("N134").Value = -12681.9132

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 657.5
$ws.Range("I113").Value = 460.3158
$ws.Range("J113").Value = 945.6923
$ws.Range("K113").Value = 1380.9474
$ws.Range("L113").Value = 2837.0769
$ws.Range("M113").Value = 789.0526
$ws.Range("N113").Value = -7177.0769
$ws.Range("H120").Value = 16506.857
$ws.Range("I120").Value = 11800
$ws.Range("J120").Value = 19121.777
$ws.Range("K120").Value = 35400
$ws.Range("L120").Value = 57365.33099999999
$ws.Range("M120").Value = -30562
$ws.Range("N120").Value = -67041.33099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 8000
$ws.Range("N70").Value = -8540
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 8000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 8000
$ws.Range("N73").Value = -9872
$ws.Range("M73").ClearContents()
$ws.Range("H126").Value = 20836322
$ws.Range("I126").Value = 83334280
$ws.Range("J126").Value = 3667.7778
$ws.Range("K126").Value = 250002840
$ws.Range("L126").Value = 11003.3334
$ws.Range("M126").Value = -250000370
$ws.Range("N126").Value = -15943.3334
$ws.Range("H132").Value = 994475.3
$ws.Range("I132").Value = 2085226.1
$ws.Range("J132").Value = 2883.6365
$ws.Range("K132").Value = 6255678.300000001
$ws.Range("L132").Value = 8650.9095
$ws.Range("M132").Value = -6253148.300000001
$ws.Range("N132").Value = -13710.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 62504064
$ws.Range("I40").Value = 83336130
$ws.Range("J40").Value = 7875
$ws.Range("K40").Value = 83336130
$ws.Range("L40").Value = 7875
$ws.Range("M40").Value = -83335994
$ws.Range("N40").Value = -8147
$ws.Range("H122").Value = 5141.4287
$ws.Range("I122").Value = 2997.5
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 8992.5
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -6542.5
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 19252872
$ws.Range("I132").Value = 34520710
$ws.Range("J132").Value = 2119.3044
$ws.Range("K132").Value = 103562130
$ws.Range("L132").Value = 6357.9132
$ws.Range("M132").Value = -103559600
$ws.Range("N132").Value = -11417.9132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 25377
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 25377
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 25377
$ws.Range("N109").Value = -28151
$ws.Range("H122").Value = 6581.391
$ws.Range("I122").Value = 4647.625
$ws.Range("J122").Value = 11001.429
$ws.Range("K122").Value = 13942.875
$ws.Range("L122").Value = 33004.287
$ws.Range("M122").Value = -11492.875
$ws.Range("N122").Value = -37904.287
$ws.Range("H132").Value = 2209.625
$ws.Range("I132").Value = 2267.9565
$ws.Range("J132").Value = 2155.96
$ws.Range("K132").Value = 6803.869499999999
$ws.Range("L132").Value = 6467.88
$ws.Range("M132").Value = -4273.869499999999
$ws.Range("N132").Value = -11527.88
$ws.Range("H136").Value = 6764952.5
$ws.Range("I136").Value = 10427851
$ws.Range("J136").Value = 2676.923
$ws.Range("K136").Value = 31283553
$ws.Range("L136").Value = 8030.768999999999
$ws.Range("M136").Value = -13130.769
$ws.Range("N136").Value = -13130.769
